$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 20 de Abril de 2020 a las 13:22"

# Update the Cataluña row (row 5) figures
$ws.Range("B5").Value = 41676
$ws.Range("C5").Value = 13934
$ws.Range("D5").Value = 23733
$ws.Range("E5").Value = 4009
